# Prod/Demo Verification Script fixes + Display Convenience fees QA results
# for all versions.
#
# For each "Generic" verification sheet, the Prod run was re-executed
# (new DateProd timestamp in column B) and a Demo result of "Fail" was
# recorded (ResultDemo in column C) for every data row that previously
# only carried a Prod result.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> rows to update, each row giving the new
# DateProd (column B) value; ResultDemo (column C) is always "Fail".
$updates = @{
    "VT-P-DebitVoid-DualCF-Generic"   = @{
        2 = "Wed May 21 15:14:57 IST 2025"
    }
    "VT-P-DebitVoid-SingleCF-Generic" = @{
        2 = "Wed May 21 15:20:23 IST 2025"
        3 = "Wed May 21 15:21:21 IST 2025"
        4 = "Wed May 21 15:22:16 IST 2025"
        5 = "Wed May 21 15:23:14 IST 2025"
    }
    "VT-P-DebitVoid-NoCF-Generic"     = @{
        2 = "Wed May 21 15:15:58 IST 2025"
        3 = "Wed May 21 15:17:02 IST 2025"
        4 = "Wed May 21 15:18:08 IST 2025"
        5 = "Wed May 21 15:19:15 IST 2025"
    }
    "VT-P-DebitCredit-DualCF-Generic" = @{
        2 = "Wed May 21 15:02:50 IST 2025"
        3 = "Wed May 21 15:03:56 IST 2025"
        4 = "Wed May 21 15:04:56 IST 2025"
        5 = "Wed May 21 15:05:59 IST 2025"
    }
    "VT-P-DebitCredit-SingleCF-Gener" = @{
        2 = "Wed May 21 15:10:52 IST 2025"
        3 = "Wed May 21 15:12:04 IST 2025"
        4 = "Wed May 21 15:13:00 IST 2025"
        5 = "Wed May 21 15:14:02 IST 2025"
    }
    "VT-P-DebitCredit-NoCF-Generic"   = @{
        2 = "Wed May 21 15:06:52 IST 2025"
        3 = "Wed May 21 15:07:51 IST 2025"
        4 = "Wed May 21 15:08:45 IST 2025"
        5 = "Wed May 21 15:09:46 IST 2025"
    }
    "VT-C-DebitCredit-DualCF-Generic" = @{
        2 = "Wed May 21 14:59:38 IST 2025"
    }
    "VT-C-DebitCredit-SingleCF-Gener" = @{
        2 = "Wed May 21 15:49:31 IST 2025"
    }
    "VT-C-DebitCredit-NoCF-Generic"   = @{
        2 = "Wed May 21 15:00:41 IST 2025"
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $dateProd = $rows[$row]
        $ws.Cells.Item($row, 2).Value = $dateProd
        $ws.Cells.Item($row, 3).Value = "Fail"
    }
}
